$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 100000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 100000
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = 100000
$ws.Range("N13").Value = -100338
$ws.Range("L13").ClearContents()

$ws.Range("H137").Value = 1916.7609
$ws.Range("I137").Value = 2188.3333
$ws.Range("J137").Value = 1620.5
$ws.Range("K137").Value = 6564.999899999999
$ws.Range("L137").Value = 4861.5
$ws.Range("M137").Value = -4014.999899999999
$ws.Range("N137").Value = -9961.5

$ws.Range("H138").Value = 2426.3674
$ws.Range("I138").Value = 1008.45
$ws.Range("J138").Value = 3404.2415
$ws.Range("K138").Value = 3025.35
$ws.Range("L138").Value = 10212.7245
$ws.Range("M138").Value = 2114.65
$ws.Range("N138").Value = -20492.7245

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -385

$ws.Range("H10").Value = 100000
$ws.Range("J10").Value = 100000
$ws.Range("L10").Value = 100000
$ws.Range("N10").Value = -100340

$ws.Range("H32").Value = 6041.7295
$ws.Range("I32").Value = 5743.3335
$ws.Range("J32").Value = 9366.714
$ws.Range("K32").Value = 5743.3335
$ws.Range("L32").Value = 9366.714
$ws.Range("M32").Value = -5456.3335
$ws.Range("N32").Value = -9940.714

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H88").Value = 3120
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3120
$ws.Range("K88").Value = 0
$ws.Range("M88").Value = 3120
$ws.Range("N88").Value = -3932
$ws.Range("L88").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H91").Value = 3120
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3120
$ws.Range("K91").Value = 0
$ws.Range("M91").Value = 3120
$ws.Range("N91").Value = -5928
$ws.Range("L91").ClearContents()

$ws.Range("H110").Value = 1757.88
$ws.Range("I110").Value = 1532.85
$ws.Range("J110").Value = 2658
$ws.Range("K110").Value = 1532.85
$ws.Range("L110").Value = 2658
$ws.Range("M110").Value = 512.1500000000001
$ws.Range("N110").Value = -6748

$ws.Range("H132").Value = 5883.2856
$ws.Range("I132").Value = 3137.8572
$ws.Range("J132").Value = 10001.429
$ws.Range("K132").Value = 9413.571599999999
$ws.Range("L132").Value = 30004.287
$ws.Range("M132").Value = -6883.571599999999
$ws.Range("N132").Value = -35064.287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1474.5385
$ws.Range("I86").Value = 1362.4445
$ws.Range("J86").Value = 1726.75
$ws.Range("K86").Value = 1362.4445
$ws.Range("L86").Value = 1726.75
$ws.Range("M86").Value = -239.4445000000001
$ws.Range("N86").Value = -3972.75

$ws.Range("H89").Value = 1474.5385
$ws.Range("I89").Value = 1362.4445
$ws.Range("J89").Value = 1726.75
$ws.Range("K89").Value = 6812.2225
$ws.Range("L89").Value = 8633.75
$ws.Range("M89").Value = -1196.2225
$ws.Range("N89").Value = -19865.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 12852.25
$ws.Range("I12").Value = 402.57144
$ws.Range("J12").Value = 100000
$ws.Range("K12").Value = 402.57144
$ws.Range("L12").Value = 100000
$ws.Range("M12").Value = -232.57144
$ws.Range("N12").Value = -100340

$ws.Range("H16").Value = 3952.7
$ws.Range("I16").Value = 5300.2
$ws.Range("J16").Value = 2605.2
$ws.Range("K16").Value = 5300.2
$ws.Range("L16").Value = 2605.2
$ws.Range("M16").Value = -5013.2
$ws.Range("N16").Value = -3179.2

$ws.Range("H62").Value = 18281.857
$ws.Range("I62").Value = 13327.333
$ws.Range("K62").Value = 13327.333
$ws.Range("M62").Value = -12703.333

$ws.Range("H65").Value = 18281.857
$ws.Range("I65").Value = 13327.333
$ws.Range("K65").Value = 66636.66500000001
$ws.Range("M65").Value = -63516.66500000001

$ws.Range("H113").Value = 3952.7
$ws.Range("I113").Value = 5300.2
$ws.Range("J113").Value = 2605.2
$ws.Range("K113").Value = 5300.2
$ws.Range("L113").Value = 2605.2
$ws.Range("M113").Value = -3130.2
$ws.Range("N113").Value = -6945.2

$ws.Range("H122").Value = 66668120
$ws.Range("I122").Value = 90909730
$ws.Range("J122").Value = 3699.75
$ws.Range("K122").Value = 272729190
$ws.Range("L122").Value = 11099.25
$ws.Range("M122").Value = -272726740
$ws.Range("N122").Value = -15999.25

$ws.Range("H132").Value = 2495.8157
$ws.Range("I132").Value = 1349.1305
$ws.Range("J132").Value = 4254.067
$ws.Range("K132").Value = 4047.3915
$ws.Range("L132").Value = 12762.201
$ws.Range("M132").Value = -1517.3915
$ws.Range("N132").Value = -17822.201

$ws.Range("H134").Value = 1472.8334
$ws.Range("I134").Value = 738.3889
$ws.Range("J134").Value = 2574.5
$ws.Range("K134").Value = 2215.1667
$ws.Range("L134").Value = 7723.5
$ws.Range("M134").Value = 319.8332999999998
$ws.Range("N134").Value = -12793.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1389.4166
$ws.Range("I50").Value = 74.71429000000001
$ws.Range("K50").Value = 224.14287
$ws.Range("M50").Value = 256.85713

$ws.Range("H53").Value = 1389.4166
$ws.Range("I53").Value = 74.71429000000001
$ws.Range("K53").Value = 224.14287
$ws.Range("M53").Value = 256.85713

$ws.Range("H98").Value = 2411
$ws.Range("I98").Value = 652.625
$ws.Range("K98").Value = 1957.875
$ws.Range("M98").Value = -459.875

$ws.Range("H113").Value = 905.7368
$ws.Range("I113").Value = 635.75
$ws.Range("J113").Value = 1368.5714
$ws.Range("K113").Value = 1907.25
$ws.Range("L113").Value = 4105.7142
$ws.Range("M113").Value = 262.75
$ws.Range("N113").Value = -8445.7142

$ws.Range("H132").Value = 2220.6
$ws.Range("I132").Value = 899.8570999999999
$ws.Range("J132").Value = 3376.25
$ws.Range("K132").Value = 8098.7139
$ws.Range("L132").Value = 30386.25
$ws.Range("M132").Value = -5568.7139
$ws.Range("N132").Value = -35446.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 6065625.5
$ws.Range("I14").Value = 6065625.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 6065625.5
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -6065457.5
$ws.Range("M14").ClearContents()

$ws.Range("H70").Value = 5801.933
$ws.Range("I70").Value = 5160
$ws.Range("J70").Value = 6292.8237
$ws.Range("K70").Value = 5160
$ws.Range("L70").Value = 6292.8237
$ws.Range("M70").Value = -4890
$ws.Range("N70").Value = -6832.8237

$ws.Range("H73").Value = 5801.933
$ws.Range("I73").Value = 5160
$ws.Range("J73").Value = 6292.8237
$ws.Range("K73").Value = 5160
$ws.Range("L73").Value = 6292.8237
$ws.Range("M73").Value = -4224
$ws.Range("N73").Value = -8164.8237

$ws.Range("H80").Value = 2976.875
$ws.Range("I80").Value = 2828.75
$ws.Range("J80").Value = 3125
$ws.Range("K80").Value = 2828.75
$ws.Range("L80").Value = 3125
$ws.Range("M80").Value = -1830.75
$ws.Range("N80").Value = -5121

$ws.Range("H83").Value = 2976.875
$ws.Range("I83").Value = 2828.75
$ws.Range("J83").Value = 3125
$ws.Range("K83").Value = 14143.75
$ws.Range("L83").Value = 15625
$ws.Range("M83").Value = -9151.75
$ws.Range("N83").Value = -25609

$ws.Range("H126").Value = 12503183
$ws.Range("I126").Value = 20837104
$ws.Range("J126").Value = 2301.75
$ws.Range("K126").Value = 62511312
$ws.Range("L126").Value = 6905.25
$ws.Range("M126").Value = -62508842
$ws.Range("N126").Value = -11845.25

$ws.Range("H132").Value = 3092.6
$ws.Range("I132").Value = 3017.3333
$ws.Range("J132").Value = 3116.3684
$ws.Range("K132").Value = 9051.999899999999
$ws.Range("L132").Value = 9349.1052
$ws.Range("M132").Value = -6521.999899999999
$ws.Range("N132").Value = -14409.1052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2349.3333
$ws.Range("I7").Value = 1945.6666
$ws.Range("K7").Value = 1945.6666
$ws.Range("M7").Value = -1833.6666

$ws.Range("H122").Value = 6336.8
$ws.Range("I122").Value = 6336.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 19010.4
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -16560.4
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 2349.3333
$ws.Range("I126").Value = 1945.6666
$ws.Range("K126").Value = 5836.9998
$ws.Range("M126").Value = -3366.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()

$ws.Range("H81").Value = 2769.8572
$ws.Range("I81").Value = 4030.3333
$ws.Range("J81").Value = 1824.5
$ws.Range("K81").Value = 8060.6666
$ws.Range("L81").Value = 3649
$ws.Range("M81").Value = -6999.6666
$ws.Range("N81").Value = -5771

$ws.Range("H84").Value = 2769.8572
$ws.Range("I84").Value = 4030.3333
$ws.Range("J84").Value = 1824.5
$ws.Range("K84").Value = 40303.333
$ws.Range("L84").Value = 18245
$ws.Range("M84").Value = -34999.333
$ws.Range("N84").Value = -28853
